# "Fix league war bug"
#
# The "war" reward columns (G = League, H = Battle) on the "league" sheet
# were paying out way too much (200 / 100, with G4:G18 even being driven by
# a stray "=F*2" shared formula). Replace them with plain, corrected
# constants: League = 25, Battle = 12 for every league row (4-18).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("league")
$ws.Activate()

# G4:G18 previously held 200 (G4 as a literal, G5:G18 as the shared formula
# "=F*2"); H4:H18 previously held 100. Overwrite with the corrected values -
# this also clears out the old formula in G5:G18, turning those cells into
# plain numbers like G4.
$ws.Range("G4:G18").Value = 25
$ws.Range("H4:H18").Value = 12

# Restore the author's cursor position / window placement at save time.
$ws.Range("S29").Select()
$excel.ActiveWindow.Left = 80
$excel.ActiveWindow.Top = 0
